$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 92.5
$ws.Range("I31").Value = 92.5
$ws.Range("K31").Value = 277.5
$ws.Range("M31").Value = -47.5
$ws.Range("H75").Value = 55000
$ws.Range("I75").Value = 55000
$ws.Range("K75").Value = 55000
$ws.Range("M75").Value = -54064
$ws.Range("H78").Value = 55000
$ws.Range("I78").Value = 55000
$ws.Range("K78").Value = 165000
$ws.Range("M78").Value = -160320
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384
$ws.Range("H96").Value = 1936
$ws.Range("I96").Value = 1936
$ws.Range("K96").Value = 5808
$ws.Range("M96").Value = -4435
$ws.Range("H135").Value = 1066.4615
$ws.Range("I135").Value = 959.9167
$ws.Range("K135").Value = 8639.2503
$ws.Range("M135").Value = -6104.2503
$ws.Range("H138").Value = 4685.185
$ws.Range("I138").Value = 3546.9
$ws.Range("J138").Value = 6108.0415
$ws.Range("K138").Value = 10640.7
$ws.Range("L138").Value = 18324.1245
$ws.Range("M138").Value = -5500.700000000001
$ws.Range("N138").Value = -28604.1245
$ws.Range("H141").Value = 2006.2
$ws.Range("I141").Value = 1970.5217
$ws.Range("J141").Value = 2416.5
$ws.Range("K141").Value = 5911.5651
$ws.Range("L141").Value = 7249.5
$ws.Range("M141").Value = -731.5650999999998
$ws.Range("N141").Value = -17609.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 799.6667
$ws.Range("I3").Value = 799.6667
$ws.Range("K3").Value = 799.6667
$ws.Range("M3").Value = -684.6667
$ws.Range("H74").Value = 1436
$ws.Range("I74").Value = 1179.6
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 1179.6
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -305.5999999999999
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 1436
$ws.Range("I77").Value = 1179.6
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 5898
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -1530
$ws.Range("N77").Value = -28736
$ws.Range("H122").Value = 7205.8237
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4702.75
$ws.Range("I20").Value = 5450.6665
$ws.Range("J20").Value = 2459
$ws.Range("K20").Value = 5450.6665
$ws.Range("L20").Value = 2459
$ws.Range("M20").Value = -5203.6665
$ws.Range("N20").Value = -2953
$ws.Range("H42").Value = 100000
$ws.Range("J42").Value = 100000
$ws.Range("L42").Value = 100000
$ws.Range("N42").Value = -100656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 449.5
$ws.Range("J22").Value = 450
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -1150

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 393.25
$ws.Range("I12").Value = 429.8
$ws.Range("K12").Value = 1289.4
$ws.Range("M12").Value = -1116.4
$ws.Range("H46").Value = 5000500
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2909
$ws.Range("H62").Value = 5665
$ws.Range("J62").Value = 4247.5
$ws.Range("L62").Value = 12742.5
$ws.Range("N62").Value = -14114.5
$ws.Range("H65").Value = 5665
$ws.Range("J65").Value = 4247.5
$ws.Range("L65").Value = 38227.5
$ws.Range("N65").Value = -45091.5
$ws.Range("H68").Value = 3276.5
$ws.Range("J68").Value = 3276.5
$ws.Range("L68").Value = 9829.5
$ws.Range("N68").Value = -11451.5
$ws.Range("H71").Value = 3276.5
$ws.Range("J71").Value = 3276.5
$ws.Range("L71").Value = 29488.5
$ws.Range("N71").Value = -37600.5
$ws.Range("H108").Value = 2262
$ws.Range("I108").Value = 1402.2222
$ws.Range("K108").Value = 4206.6666
$ws.Range("M108").Value = -1326.6666
$ws.Range("H113").Value = 3341.3333
$ws.Range("J113").Value = 3029.8
$ws.Range("L113").Value = 9089.400000000001
$ws.Range("N113").Value = -13429.4
$ws.Range("H117").Value = 1229.0834
$ws.Range("J117").Value = 1677.5
$ws.Range("L117").Value = 5032.5
$ws.Range("N117").Value = -11916.5
$ws.Range("H126").Value = 4289.8
$ws.Range("I126").Value = 3362.25
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 10086.75
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -5146.75
$ws.Range("N126").Value = -33880
$ws.Range("H129").Value = 1359.1666
$ws.Range("J129").Value = 1856.75
$ws.Range("L129").Value = 5570.25
$ws.Range("N129").Value = -15570.25
$ws.Range("H130").Value = 5899.4
$ws.Range("I130").Value = 999
$ws.Range("J130").Value = 9166.333000000001
$ws.Range("K130").Value = 2997
$ws.Range("L130").Value = 27498.999
$ws.Range("M130").Value = 2023
$ws.Range("N130").Value = -37538.999
$ws.Range("H131").Value = 1445.8
$ws.Range("I131").Value = 650
$ws.Range("K131").Value = 1950
$ws.Range("M131").Value = 3090
$ws.Range("H134").Value = 1121.1666
$ws.Range("I134").Value = 1121.1666
$ws.Range("K134").Value = 3363.4998
$ws.Range("M134").Value = 1706.5002
$ws.Range("H137").Value = 6400
$ws.Range("I137").Value = 6666.6665
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 19999.9995
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -14899.9995
$ws.Range("N137").Value = -28200
$ws.Range("H138").Value = 4099.5
$ws.Range("I138").Value = 4099.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 12298.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -7158.5
$ws.Range("N138").Value = ""
$ws.Range("H140").Value = 1314.5555
$ws.Range("I140").Value = 1314.5555
$ws.Range("K140").Value = 3943.6665
$ws.Range("M140").Value = 1236.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 63.333332
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 60
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = 52
$ws.Range("N5").Value = -289
$ws.Range("H33").Value = 10285.143
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 10285.143
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 10285.143
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = -10789.143
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = ""
$ws.Range("H70").Value = 6292
$ws.Range("J70").Value = 7333.1665
$ws.Range("L70").Value = 7333.1665
$ws.Range("N70").Value = -7873.1665
$ws.Range("H73").Value = 6292
$ws.Range("J73").Value = 7333.1665
$ws.Range("L73").Value = 7333.1665
$ws.Range("N73").Value = -9205.166499999999
$ws.Range("H80").Value = 16766.445
$ws.Range("I80").Value = 7374.75
$ws.Range("J80").Value = 24279.8
$ws.Range("K80").Value = 7374.75
$ws.Range("L80").Value = 24279.8
$ws.Range("M80").Value = -6376.75
$ws.Range("N80").Value = -26275.8
$ws.Range("H83").Value = 16766.445
$ws.Range("I83").Value = 7374.75
$ws.Range("J83").Value = 24279.8
$ws.Range("K83").Value = 36873.75
$ws.Range("L83").Value = 121399
$ws.Range("M83").Value = -31881.75
$ws.Range("N83").Value = -131383
$ws.Range("H122").Value = 66996.94
$ws.Range("I122").Value = 3698.3
$ws.Range("K122").Value = 11094.9
$ws.Range("M122").Value = -8644.900000000001
$ws.Range("H132").Value = 2169.9285
$ws.Range("I132").Value = 2130.4167
$ws.Range("J132").Value = 2407
$ws.Range("K132").Value = 6391.250100000001
$ws.Range("L132").Value = 7221
$ws.Range("M132").Value = -3861.250100000001
$ws.Range("N132").Value = -12281

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2705.575
$ws.Range("I22").Value = 1647.5625
$ws.Range("J22").Value = 6937.625
$ws.Range("K22").Value = 1647.5625
$ws.Range("L22").Value = 6937.625
$ws.Range("M22").Value = -1352.5625
$ws.Range("N22").Value = -7527.625
$ws.Range("H25").Value = 13885
$ws.Range("I25").Value = 9862.5
$ws.Range("J25").Value = 16566.666
$ws.Range("K25").Value = 9862.5
$ws.Range("L25").Value = 16566.666
$ws.Range("M25").Value = -9632.5
$ws.Range("N25").Value = -17026.666
$ws.Range("H27").Value = 2705.575
$ws.Range("I27").Value = 1647.5625
$ws.Range("J27").Value = 6937.625
$ws.Range("K27").Value = 1647.5625
$ws.Range("L27").Value = 6937.625
$ws.Range("M27").Value = -1540.5625
$ws.Range("N27").Value = -7151.625
$ws.Range("H46").Value = 4159.7
$ws.Range("I46").Value = 3108.3333
$ws.Range("J46").Value = 4610.2856
$ws.Range("K46").Value = 3108.3333
$ws.Range("L46").Value = 4610.2856
$ws.Range("M46").Value = -2920.3333
$ws.Range("N46").Value = -4986.2856
$ws.Range("H82").Value = 3034.25
$ws.Range("I82").Value = 3795.6667
$ws.Range("J82").Value = 750
$ws.Range("K82").Value = 3795.6667
$ws.Range("L82").Value = 750
$ws.Range("M82").Value = -3434.6667
$ws.Range("N82").Value = -1472
$ws.Range("H85").Value = 3034.25
$ws.Range("I85").Value = 3795.6667
$ws.Range("J85").Value = 750
$ws.Range("K85").Value = 3795.6667
$ws.Range("L85").Value = 750
$ws.Range("M85").Value = -2547.6667
$ws.Range("N85").Value = -3246
$ws.Range("H122").Value = 10666
$ws.Range("I122").Value = 10666
$ws.Range("K122").Value = 31998
$ws.Range("M122").Value = -29548

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 33352.332
$ws.Range("J28").Value = 33352.332
$ws.Range("L28").Value = 33352.332
$ws.Range("N28").Value = -34048.332
